$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, matching the formatting of the other header
# cells in row 1 (e.g. G1 "sum") by copying that cell's format over.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add numeric value 0 in H2 (plain, unstyled like the other data cells)
$ws.Range("H2").Value = 0
